$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 607 (weekly price update), shifting the
# existing rows 607-646 down to 608-647.
$ws.Rows.Item(607).Insert()

# Populate the new row 607 with data for the new week (2023-12-05),
# a copy of the prior "Primera" quality entry with updated price.
$ws.Range("A607").Value = 7
$ws.Range("B607").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C607").Value = "Ñuble"
$ws.Range("D607").Value = 45265
$ws.Range("E607").Value = 16
$ws.Range("F607").Value = 100112009
$ws.Range("G607").Value = "Acelga"
$ws.Range("H607").Value = "Sin especificar"
$ws.Range("I607").Value = "Primera"
$ws.Range("J607").Value = 200
$ws.Range("K607").Value = 700
$ws.Range("L607").Value = 700
$ws.Range("M607").Value = 700
$ws.Range("N607").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O607").Value = "Región de Ñuble"
$ws.Range("P607").Value = 700
$ws.Range("Q607").Value = 1
$ws.Range("R607").Value = "Hortaliza"
